# NYPD 62nd Precinct weekly CompStat report — roll the report forward one
# week (Volume 30, Number 27 -> 28; week of 7/3/2023-7/9/2023 ->
# 7/10/2023-7/16/2023) and refresh the crime-complaint figures that came
# in with the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: issue number and reporting week -------------------------------
$ws.Range("A8").Value = "Volume 30   Number  28"
$ws.Range("C9").Value = "Report Covering the Week  7/10/2023  Through  7/16/2023"

# --- Cells that flip between "no activity" text (***.*) and a real number --
# Donor cells supply the number format / style so the converted cell keeps
# the workbook's existing look (General-text style 14 for text, the
# thousands style for whole numbers, the decimal style for percentages).
$numWholeDonor = $ws.Range("C16")   # style used for whole-number counts
$numPctDonor   = $ws.Range("H26")   # style used for signed decimal percentages
$textZeroDonor = $ws.Range("D22")   # style used for the literal "0" text
$textNaDonor   = $ws.Range("E30")   # style used for the literal "***.*" text

# Set the value first (while a plain/unambiguous number format is active
# so Excel stores it as a true number or true text, not a look-alike of
# the other type), *then* paste-special just the formatting from a donor
# cell that already has the style we want — this swaps the style index
# without disturbing the freshly-written value/type.
function Convert-ToNumber($range, $donor, $value) {
    $range.NumberFormat = "General"
    $range.Value = $value
    $donor.Copy() | Out-Null
    $range.PasteSpecial(-4122) | Out-Null
}

function Convert-ToText($range, $donor, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $donor.Copy() | Out-Null
    $range.PasteSpecial(-4122) | Out-Null
}

# Row 22 (Transit): weekly 2023 count goes from a reported 1 down to none.
Convert-ToText $ws.Range("C22") $textZeroDonor "0"

# Row 26 (UCR Rape*): weekly 2022 count now has a real (non-zero) figure.
Convert-ToNumber $ws.Range("D26") $numWholeDonor 1
Convert-ToNumber $ws.Range("E26") $numPctDonor -100

# Row 27 (Other Sex Crimes): weekly 2023 count now has a real figure.
Convert-ToNumber $ws.Range("C27") $numWholeDonor 2

# Row 28 (Shooting Vic.): weekly 2022 count drops back to none.
Convert-ToText $ws.Range("D28") $textZeroDonor "0"
Convert-ToText $ws.Range("E28") $textNaDonor "***.*"

# Row 29 (Shooting Inc.): weekly 2022 count drops back to none.
Convert-ToText $ws.Range("D29") $textZeroDonor "0"
Convert-ToText $ws.Range("E29") $textNaDonor "***.*"

# --- Plain numeric refreshes (value only, style unchanged) -----------------
$values = @{
    "G15" = 1

    "C16" = 3;  "D16" = 1;   "E16" = 200;                 "F16" = 9;   "G16" = 11
    "H16" = -18.181818181818; "I16" = 75;  "J16" = 62
    "K16" = 20.967741935483;  "L16" = 36.363636363636
    "M16" = -14.772727272727; "N16" = -82.517482517482

    "C17" = 5;  "D17" = 7;   "E17" = -28.571428571428;    "F17" = 17;  "G17" = 22
    "H17" = -22.727272727272; "I17" = 116; "J17" = 93
    "K17" = 24.731182795698;  "L17" = 14.851485148514
    "M17" = 90.163934426229;  "N17" = -29.268292682926

    "C18" = 6;  "D18" = 5;   "E18" = 20;                  "F18" = 12;  "G18" = 7
    "H18" = 71.428571428571;  "I18" = 122; "J18" = 93
    "K18" = 31.182795698924;  "L18" = 56.410256410256
    "M18" = -22.292993630573; "N18" = -84.398976982097

    "C19" = 13; "D19" = 10;  "E19" = 30;                  "F19" = 51;  "G19" = 49
    "H19" = 4.081632653061;   "I19" = 350; "J19" = 375
    "K19" = -6.666666666666;  "L19" = 18.243243243243
    "M19" = 54.185022026431;  "N19" = -13.580246913580

    "C20" = 8;  "D20" = 3;   "E20" = 166.666666666667;    "F20" = 16;  "G20" = 17
    "H20" = -5.882352941176;  "I20" = 80;  "J20" = 74
    "K20" = 8.108108108108;   "L20" = 100
    "M20" = -13.978494623655; "N20" = -91.726990692864

    "C21" = 35; "D21" = 26;  "E21" = 34.615384615384;     "F21" = 105; "G21" = 107
    "H21" = -1.869158878504;  "I21" = 752; "J21" = 711
    "K21" = 5.766526019690;   "L21" = 30.103806228373
    "M21" = 18.238993710691;  "N21" = -72.773352643012

    "F22" = 2;  "H22" = 100

    "C24" = 27; "D24" = 40;  "E24" = -32.5;                "F24" = 136; "G24" = 178
    "H24" = -23.595505617977; "I24" = 975; "J24" = 926
    "K24" = 5.291576673866;   "L24" = 51.632970451010
    "M24" = 95.783132530120

    "C25" = 7;                "E25" = -30;                 "F25" = 44;  "G25" = 35
    "H25" = 25.714285714285;  "I25" = 324; "J25" = 238
    "K25" = 36.134453781512;  "L25" = 38.461538461538
    "M25" = 26.5625

    "J26" = 14; "K26" = -21.428571428571

    "F27" = 10; "I27" = 29;  "K27" = 11.538461538461;     "L27" = -6.451612903225
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
